$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# Row 31
$ws.Range("A31").Value = 41466
$ws.Range("A31").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B31").Value = 2.5
$ws.Range("C31").Value = 1.5
$ws.Range("D31").Value = "Implementation tc14"

# Row 32
$ws.Range("A32").Value = 41467
$ws.Range("A32").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("C32").Value = 4.25
$ws.Range("D32").Value = "Implementation tc14"

$ws.Range("C32").Select()
